$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.047, 0.026, 0.522, 99, -33, 48),
    @(-0.479, -0.141, 0.687, 46, -62, 72),
    @(-0.062, -0.307, 0.759, 110, -11, 72),
    @(0.043, 0.153, 0.496, 107, 74, -140),
    @(0.092, 0.008, 0.281, 126, 0, -2),
    @(0.378, 0.208, 0.611, -32, 76, 123),
    @(-0.523, 0.197, -0, -150, 11, -170),
    @(-0.354, -0.251, 0.614, 69, -5, 29),
    @(-0.223, 0.194, 0.067, -50, 0, 34),
    @(-0.285, -0.168, 0.683, -60, -13, 169)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
